$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new value (kept as text, same as original inline strings)
$updates = @{
    "D2"  = "303.24"
    "D3"  = "31.72"
    "E3"  = "0.36%"
    "D4"  = "5.164"
    "E4"  = "0.61%"
    "D5"  = "0.07825"
    "E5"  = "4.23%"
    "D6"  = "2.328"
    "E6"  = "36.93%"
    "D7"  = "7.951"
    "E7"  = "2.76%"
    "D8"  = "3.872"
    "E8"  = "1.86%"
    "D9"  = "0.9087"
    "E9"  = "-2.49%"
    "D10" = "0.1730"
    "E10" = "2.27%"
    "D11" = "0.07324"
    "E11" = "0.93%"
    "D12" = "0.08123"
    "E12" = "2.54%"
    "D13" = "0.03024"
    "E13" = "0.26%"
    "D14" = "0.09947"
    "E14" = "0.44%"
    "D15" = "0.001506"
    "E15" = "0.38%"
    "D16" = "0.006024"
    "E16" = "-5.46%"
    "D17" = "3.497"
    "E17" = "1.50%"
    "D18" = "2.238"
    "E18" = "0.54%"
    "D19" = "0.3244"
    "E19" = "-1.23%"
    "E20" = "0.79%"
    "D21" = "4.662"
    "E21" = "1.87%"
    "D22" = "0.04652"
    "E22" = "-0.24%"
    "E23" = "0.36%"
    "D24" = "0.001259"
    "E24" = "3.40%"
    "D25" = "0.004518"
    "E25" = "2.17%"
    "E26" = "3.63%"
    "D27" = "0.0002741"
    "D39" = "0.01790"
    "E39" = "7.02%"
    "D40" = "0.04561"
    "E40" = "2.37%"
    "D41" = "0.007271"
    "E41" = "2.91%"
    "D42" = "0.1361"
    "E42" = "2.65%"
    "D43" = "0.002237"
    "E43" = "8.45%"
    "D44" = "0.01075"
    "E44" = "-4.72%"
    "D45" = "0.00006499"
    "E45" = "8.00%"
    "D46" = "0.00000000750"
    "E46" = "-0.06%"
    "E47" = "-57.22%"
    "E49" = "-0.06%"
    "E50" = "0.01%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force the cell to Text format first so Excel does not reinterpret the
    # numeric/percentage-looking string as a real number, preserving the
    # exact text content (matches the original inline string representation).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
